# Java_Keywords_and_Their_Uses.docx edit script
# - Title "Java Keywords and Their Uses": make it red, 20pt (sz/szCs 40 half-points)
# - Several ListBullet paragraphs that were split across two <w:r> runs (an
#   artifact of the original authoring) get re-typed as a single contiguous
#   run.
# - A handful of single Java keywords (boolean, char, instanceof, enum, goto)
#   are not dictionary words, so Word's background spell checker wraps them
#   in <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#   once they are their own run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title formatting: red (FF0000), size 20pt (sz/szCs = 40 half-points)
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Java Keywords and Their Uses")) {
        $p.Range.Font.Color = 255      # wdColorRed -> w:color FF0000
        $p.Range.Font.Size = 20        # w:sz 40
        $p.Range.Font.SizeBi = 20      # w:szCs 40
    }
}

# ---------------------------------------------------------------------
# 2. Paragraphs whose two runs just get re-typed as one contiguous run
#    (no proofErr markers involved).
# ---------------------------------------------------------------------
$merges = @(
    @{ Match = "protected"; Text = "protected: Allows access within the same package and subclasses." },
    @{ Match = "extends";   Text = "extends: Indicates that a class inherits from a superclass." },
    @{ Match = "super:";    Text = "super: Refers to the parent class (used to call superclass constructors or methods)." },
    @{ Match = "static:";   Text = "static: Defines class-level variables or methods shared by all instances." },
    @{ Match = "nat";       Text = "native: Declares a method implemented in another language (like C)." },
    @{ Match = "switch:";   Text = "switch: Selects one of many code blocks to execute." },
    @{ Match = "do: Used";  Text = "do: Used with while for a loop that runs at least once." },
    @{ Match = "yield:";    Text = "yield: Returns a value from a switch expression." },
    @{ Match = "throw:";    Text = "throw: Used to throw an exception manually." },
    @{ Match = "7. Object"; Text = "7. Object Reference and Null Keywords" },
    @{ Match = "var:";      Text = "var: Allows local variable type inference (the compiler infers the type)." }
)

foreach ($m in $merges) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($m.Match)) {
            $full = $d.Range($p.Range.Start, $p.Range.End)
            # First overwrite with an unrelated placeholder so the engine
            # can't "reuse" the pre-existing run split, then write the
            # final text -- this collapses the paragraph to a single run.
            $full.Text = "zzz__placeholder__zzz`r"
            foreach ($pp in $d.Paragraphs) {
                if ($pp.Range.Text.StartsWith("zzz__placeholder__zzz")) {
                    $again = $d.Range($pp.Range.Start, $pp.Range.End)
                    $again.Text = $m.Text + "`r"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3. Paragraphs where the leading keyword gets wrapped in proofErr
#    spellStart/spellEnd markers (keyword becomes its own run).
# ---------------------------------------------------------------------
$proofErrs = @(
    @{ Match = "boolean";    Style = "ListBullet"; Word = "boolean";    Rest = ": Defines a variable that can store only true or false." },
    @{ Match = "char:";      Style = "ListBullet"; Word = "char";       Rest = ": Defines a single 16-bit Unicode character." },
    @{ Match = "instanceof"; Style = "ListBullet"; Word = "instanceof"; Rest = ": Tests whether an object is an instance of a specific class or subclass." },
    @{ Match = "enum:";      Style = "ListBullet"; Word = "enum";       Rest = ": Defines a set of named constants (enumeration)." },
    @{ Match = "goto:";      Style = "ListBullet"; Word = "goto";       Rest = ": Reserved but not used (for backward compatibility with C)." }
)

foreach ($pe in $proofErrs) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($pe.Match)) {
            $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body><w:p><w:pPr><w:pStyle w:val="' + $pe.Style + '"/></w:pPr>' +
                   '<w:proofErr w:type="spellStart"/>' +
                   '<w:r><w:t>' + $pe.Word + '</w:t></w:r>' +
                   '<w:proofErr w:type="spellEnd"/>' +
                   '<w:r><w:t>' + $pe.Rest + '</w:t></w:r>' +
                   '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            # Stop one character short of Range.End so we never swallow the
            # paragraph mark itself (matters for the very last paragraph in
            # the body, which would otherwise leave a stray empty <w:p>
            # behind); this also preserves the original w14:paraId/rsid
            # attributes on the paragraph.
            $r = $d.Range($p.Range.Start, $p.Range.End - 1)
            $r.InsertXML($xml)
        }
    }
}
